# Add a new "HeroSkillLevel" column (E) to the LevelExp table, mirroring
# the existing Id/Exp/CardExp/TowerLevel layout:
#   E1 = header "HeroSkillLevel"
#   E2 = "int"               (type row, styled like D2)
#   E3 = "英雄技能等级"       (description row, styled like D3)
#   E4:E102 = per-level hero skill level data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -----------------------------------------------------------
$ws.Range("E1").Value = "HeroSkillLevel"

# --- type / description rows, copying formatting from column D ------------
$ws.Range("E2").Value = "int"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("E3").Value = "英雄技能等级"
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- data rows (E4:E102) ---------------------------------------------------
$heroSkillLevels = @(1,1,1,1,1,2,2,2,2,2,2,2,2,3,3,3,3,3,3,4,4,4,4,4,4,4,4,4,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5)

for ($i = 0; $i -lt $heroSkillLevels.Length; $i++) {
    $ws.Cells.Item(4 + $i, 5).Value = $heroSkillLevels[$i]
}

# --- extend the table "表1" to include the new column ----------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E102"))
$ws.Range("E1").Value = "HeroSkillLevel"

# --- move the selection to where the author ended up -----------------------
$ws.Range("E31").Select()
